$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Class")

# --- New header columns: G=BatchName, H=ClassTopic, I=StaffName --------
$ws.Range("G1").Value = "BatchName"
$ws.Range("H1").Value = "ClassTopic"
$ws.Range("I1").Value = "StaffName"

# Match the formatting of the existing header cells (B1:F1) for the new
# header cells.
$ws.Range("F1").Copy()
$ws.Range("G1:I1").PasteSpecial(-4122)

# --- New scenario rows ---------------------------------------------------
# Row 7: Search with valid batch name
$ws.Cells.Item(7, 1).Value = "Search with valid batch name"
$ws.Cells.Item(7, 7).Value = "Micro service-01"
$ws.Cells.Item(7, 8).Value = "Vidhya Test"
$ws.Cells.Item(7, 9).Value = "Getha Takur"

# Row 8: Search with valid class topic
$ws.Cells.Item(8, 1).Value = "Search with valid class topic"
$ws.Cells.Item(8, 7).Value = "Micro service-01"
$ws.Cells.Item(8, 8).Value = "Vidhya Test"
$ws.Cells.Item(8, 9).Value = "Getha Takur"

# Row 9: Search with valid staff name
$ws.Cells.Item(9, 1).Value = "Search with valid staff name"
$ws.Cells.Item(9, 7).Value = "Micro service-01"
$ws.Cells.Item(9, 8).Value = "Vidhya Test"
$ws.Cells.Item(9, 9).Value = "Getha Takur"

# Match formatting: column A scenario-name cells use the same style as the
# existing scenario rows; G:I data cells use the same style as the header /
# other data cells in those columns.
$ws.Range("A2").Copy()
$ws.Range("A7:A9").PasteSpecial(-4122)

$ws.Range("G1").Copy()
$ws.Range("G7:I9").PasteSpecial(-4122)
